# edit.ps1 - apply the PtX_demand_DK Sheet1 restructuring described by the diff:
#   * insert a "Fossil Gases" row after every "Biogenic Gases" row (one per year block)
#   * insert a "Fossil Liquids" row after every "Biogenic Liquids" row (one per year block)
#   * refresh every data row (new Fossil-row figures + updated Aviation/Overall Demand totals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($sheet, $rowIndex, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $colIndex = $i + 1
        $cellValue = $values[$i]
        if ($cellValue -ne $null) {
            $sheet.Cells.Item($rowIndex, $colIndex).Value = $cellValue
        }
    }
}

# Insert the 6 new rows at their final target row numbers, top-to-bottom, so each
# index already reflects the shift caused by the inserts processed before it.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(34).Insert()

# Write final values for every data row (2-37); columns left $null stay blank.
Set-RowData $ws 2 @("Hydrogen", 2030, $null, $null, $null, 0.0004249492896250961, $null, 0.000000001306336170096178, 0.0001748509948835993, $null, $null)  # Hydrogen 2030
Set-RowData $ws 3 @("Methanol", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Methanol 2030
Set-RowData $ws 4 @("Ammonia", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Ammonia 2030
Set-RowData $ws 5 @("Synthetic Gases", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Synthetic Gases 2030
Set-RowData $ws 6 @("Biogenic Gases", 2030, $null, $null, 0.0003844825721781027, 0.0001111623516295486, $null, $null, 0.00003373774600869086, $null, $null)  # Biogenic Gases 2030
Set-RowData $ws 7 @("Fossil Gases", 2030, $null, $null, $null, 0.001635641213188201, $null, $null, 0.0001334526910548289, $null, $null)  # Fossil Gases 2030
Set-RowData $ws 8 @("Synthetic Liquids", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Synthetic Liquids 2030
Set-RowData $ws 9 @("Biogenic Liquids", 2030, $null, $null, $null, 0.007470683708903277, 0.00007937091197937183, 0.003168171690874954, 0.0062820907408898, 0.00005039905455642576, 0.004500659347109801)  # Biogenic Liquids 2030
Set-RowData $ws 10 @("Fossil Liquids", 2030, $null, $null, $null, 0.07584242857359601, 0.0005589860579383, 0.0288833029545226, 0.0398321786500645, 0.0003050655832605, 0.0437396477109587)  # Fossil Liquids 2030
Set-RowData $ws 11 @("Biomass [Solid]", 2030, $null, $null, 0.001572628624304701, $null, $null, $null, $null, $null, $null)  # Biomass [Solid] 2030
Set-RowData $ws 12 @("Renewable Energy Carrier", 2030, $null, $null, 0.001047081741485413, $null, $null, $null, $null, $null, $null)  # Renewable Energy Carrier 2030
Set-RowData $ws 13 @("Overall Demand", 2030, $null, $null, 0.003004192937968216, 0.08548486513694213, 0.0006383569699176719, 0.03205147595173372, 0.04645631082290141, 0.0003554646378169257, 0.0482403070580685)  # Overall Demand 2030
Set-RowData $ws 14 @("Hydrogen", 2040, $null, $null, $null, 0.002053604010160708, $null, 0.0000001093547660215742, 0.0002577736165919, $null, $null)  # Hydrogen 2040
Set-RowData $ws 15 @("Methanol", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Methanol 2040
Set-RowData $ws 16 @("Ammonia", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Ammonia 2040
Set-RowData $ws 17 @("Synthetic Gases", 2040, $null, $null, $null, 0.0000000007322806083117311, $null, $null, 0.00000000008152913431792302, $null, $null)  # Synthetic Gases 2040
Set-RowData $ws 18 @("Biogenic Gases", 2040, $null, $null, 0.001497759788412087, 0.0001417368229084807, $null, $null, 0.00005615015649739315, $null, $null)  # Biogenic Gases 2040
Set-RowData $ws 19 @("Fossil Gases", 2040, $null, $null, $null, 0.0008703050421561751, $null, $null, 0.0001420473188329681, $null, $null)  # Fossil Gases 2040
Set-RowData $ws 20 @("Synthetic Liquids", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Synthetic Liquids 2040
Set-RowData $ws 21 @("Biogenic Liquids", 2040, $null, $null, $null, 0.003109201685836332, 0.000129387680729, 0.0038806073302648, 0.0041737399182243, 0.00006122667160201211, 0.005115832688583799)  # Biogenic Liquids 2040
Set-RowData $ws 22 @("Fossil Liquids", 2040, $null, $null, $null, 0.020121667103184, 0.0006009410903166, 0.027242133265133, 0.0179156096703749, 0.000270718412435, 0.0423519848862636)  # Fossil Liquids 2040
Set-RowData $ws 23 @("Biomass [Solid]", 2040, $null, $null, 0.001564790190102735, $null, $null, $null, $null, $null, $null)  # Biomass [Solid] 2040
Set-RowData $ws 24 @("Renewable Energy Carrier", 2040, $null, $null, 0.00407055340951407, $null, $null, $null, $null, $null, $null)  # Renewable Energy Carrier 2040
Set-RowData $ws 25 @("Overall Demand", 2040, $null, $null, 0.007133103388028891, 0.02629651539652631, 0.0007303287710456, 0.03112284995016382, 0.0225453207620506, 0.0003319450840370121, 0.0474678175748474)  # Overall Demand 2040
Set-RowData $ws 26 @("Hydrogen", 2050, $null, $null, $null, 0.0028429310518253, $null, 0.0000001853475200557468, 0.0004117404118114, $null, $null)  # Hydrogen 2050
Set-RowData $ws 27 @("Methanol", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Methanol 2050
Set-RowData $ws 28 @("Ammonia", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null)  # Ammonia 2050
Set-RowData $ws 29 @("Synthetic Gases", 2050, $null, $null, $null, 0.000000005159994887463009, $null, $null, 0.000000002181042364749701, $null, $null)  # Synthetic Gases 2050
Set-RowData $ws 30 @("Biogenic Gases", 2050, $null, $null, 0.003608512460127723, 0.00002200879955645684, $null, $null, 0.00001599038756056599, $null, $null)  # Biogenic Gases 2050
Set-RowData $ws 31 @("Fossil Gases", 2050, $null, $null, $null, 0.00004522630829141753, $null, $null, 0.00005229533163129361, $null, $null)  # Fossil Gases 2050
Set-RowData $ws 32 @("Synthetic Liquids", 2050, $null, $null, $null, 0.00000000001721647041242079, 0.000000000004764648946044249, 0.0000000001475616247131327, 0.0000000000685489413995898, 0.0000000000004193308022998501, 0.0000000003671568424888381)  # Synthetic Liquids 2050
Set-RowData $ws 33 @("Biogenic Liquids", 2050, $null, $null, $null, 0.000273505592033718, 0.0002311930686831, 0.0051125507058801, 0.0010779854908764, 0.00007874719523598319, 0.0072768487878193)  # Biogenic Liquids 2050
Set-RowData $ws 34 @("Fossil Liquids", 2050, $null, $null, $null, 0.0010985726600198, 0.0005420234133979, 0.0246011132586586, 0.0031809784935991, 0.0002330211164886, 0.0393920586321853)  # Fossil Liquids 2050
Set-RowData $ws 35 @("Biomass [Solid]", 2050, $null, $null, 0.001549435461082234, $null, $null, $null, $null, $null, $null)  # Biomass [Solid] 2050
Set-RowData $ws 36 @("Renewable Energy Carrier", 2050, $null, $null, 0.009801744138411114, $null, $null, $null, $null, $null, $null)  # Renewable Energy Carrier 2050
Set-RowData $ws 37 @("Overall Demand", 2050, $null, $null, 0.01495969205962107, 0.00428224958893805, 0.0007732164868456489, 0.02971384945962038, 0.004738992365070066, 0.000311768312143914, 0.04666890778716144)  # Overall Demand 2050
